# Rename "Tabelle1" -> "Time Windows Comparison" and fix up the chart's
# series formulas so the cached references keep pointing at the renamed
# sheet (Excel normally does this automatically on a sheet rename; the
# chart's SERIES() formulas are updated explicitly here for parity).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$oldName = $ws.Name
$newName = "Time Windows Comparison"
$ws.Name = $newName

$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

foreach ($ser in $chart.SeriesCollection()) {
    $ser.Formula = $ser.Formula.Replace($oldName + "!", "'" + $newName + "'!")
}
